$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2079.9
$ws.Range("I9").Value = 2255.4443
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 2255.4443
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -2086.4443
$ws.Range("N9").Value = -838

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 404.92856
$ws.Range("I17").Value = 63
$ws.Range("J17").Value = 417.5926
$ws.Range("K17").Value = 189
$ws.Range("L17").Value = 1252.7778
$ws.Range("M17").Value = -21
$ws.Range("N17").Value = -1588.7778

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 13399.4
$ws.Range("I32").Value = 8186.6
$ws.Range("J32").Value = 16005.8
$ws.Range("K32").Value = 8186.6
$ws.Range("L32").Value = 16005.8
$ws.Range("M32").Value = -7860.6
$ws.Range("N32").Value = -16657.8

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1640.8182
$ws.Range("I43").Value = 1418.5
$ws.Range("J43").Value = 1690.2222
$ws.Range("K43").Value = 1418.5
$ws.Range("L43").Value = 1690.2222
$ws.Range("M43").Value = -1349.5
$ws.Range("N43").Value = -1828.2222

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 12919.5
$ws.Range("I51").Value = 5623.75
$ws.Range("J51").Value = 17783.334
$ws.Range("K51").Value = 5623.75
$ws.Range("L51").Value = 17783.334
$ws.Range("M51").Value = -5139.75
$ws.Range("N51").Value = -18751.334

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5281936
$ws.Range("I62").Value = 16998398
$ws.Range("J62").Value = 9528.200000000001
$ws.Range("K62").Value = 16998398
$ws.Range("L62").Value = 9528.200000000001
$ws.Range("M62").Value = -16997774
$ws.Range("N62").Value = -10776.2

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5281936
$ws.Range("I65").Value = 16998398
$ws.Range("J65").Value = 9528.200000000001
$ws.Range("K65").Value = 84991990
$ws.Range("L65").Value = 47641
$ws.Range("M65").Value = -84988870
$ws.Range("N65").Value = -53881

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 17869248
$ws.Range("I74").Value = 28579796
$ws.Range("J74").Value = 18333.334
$ws.Range("K74").Value = 28579796
$ws.Range("L74").Value = 18333.334
$ws.Range("M74").Value = -28578860
$ws.Range("N74").Value = -20205.334

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 17869248
$ws.Range("I77").Value = 28579796
$ws.Range("J77").Value = 18333.334
$ws.Range("K77").Value = 142898980
$ws.Range("L77").Value = 91666.67
$ws.Range("M77").Value = -142894300
$ws.Range("N77").Value = -101026.67

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3732.6875
$ws.Range("I98").Value = 3687.9167
$ws.Range("J98").Value = 3867
$ws.Range("K98").Value = 3687.9167
$ws.Range("L98").Value = 3867
$ws.Range("M98").Value = -2189.9167
$ws.Range("N98").Value = -6863

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3732.6875
$ws.Range("I122").Value = 3687.9167
$ws.Range("J122").Value = 3867
$ws.Range("K122").Value = 11063.7501
$ws.Range("L122").Value = 11601
$ws.Range("M122").Value = -8613.750100000001
$ws.Range("N122").Value = -16501

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 902.8570999999999
$ws.Range("I131").Value = 535.8333
$ws.Range("J131").Value = 3105
$ws.Range("K131").Value = 1607.4999
$ws.Range("L131").Value = 9315
$ws.Range("M131").Value = 3432.5001
$ws.Range("N131").Value = -19395

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 502667.88
$ws.Range("I132").Value = 648568.25
$ws.Range("J132").Value = 16333.333
$ws.Range("K132").Value = 1945704.75
$ws.Range("L132").Value = 48999.999
$ws.Range("M132").Value = -1943174.75
$ws.Range("N132").Value = -54059.999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3496.228
$ws.Range("I138").Value = 1892.625
$ws.Range("J138").Value = 4122.0244
$ws.Range("K138").Value = 5677.875
$ws.Range("L138").Value = 12366.0732
$ws.Range("M138").Value = -537.875
$ws.Range("N138").Value = -22646.0732

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 147.4
$ws.Range("I4").Value = 158.125
$ws.Range("J4").Value = 104.5
$ws.Range("K4").Value = 158.125
$ws.Range("L4").Value = 104.5
$ws.Range("M4").Value = -42.125
$ws.Range("N4").Value = -336.5

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 11647.571
$ws.Range("I46").Value = 22222
$ws.Range("J46").Value = 9885.166999999999
$ws.Range("K46").Value = 22222
$ws.Range("L46").Value = 9885.166999999999
$ws.Range("M46").Value = -21903
$ws.Range("N46").Value = -10523.167

# ARM row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1600522.6
$ws.Range("I134").Value = 2277003.8
$ws.Range("J134").Value = 22066.666
$ws.Range("K134").Value = 6831011.399999999
$ws.Range("L134").Value = 66199.99800000001
$ws.Range("M134").Value = -6828476.399999999
$ws.Range("N134").Value = -71269.99800000001

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3474.3333
$ws.Range("I31").Value = 1077.5652
$ws.Range("J31").Value = 5443.107
$ws.Range("K31").Value = 1077.5652
$ws.Range("L31").Value = 5443.107
$ws.Range("M31").Value = -782.5652
$ws.Range("N31").Value = -6033.107

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3474.3333
$ws.Range("I34").Value = 1077.5652
$ws.Range("J34").Value = 5443.107
$ws.Range("K34").Value = 1077.5652
$ws.Range("L34").Value = 5443.107
$ws.Range("M34").Value = -875.5652
$ws.Range("N34").Value = -5847.107

# CRP row 63
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

# CRP row 66
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 162851.58
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 162851.58
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 488554.74
$ws.Range("N37").Value = -488778.74

# CUL row 88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 7679.4
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 9224.25
$ws.Range("K88").Value = 4500
$ws.Range("L88").Value = 27672.75
$ws.Range("M88").Value = -4072
$ws.Range("N88").Value = -28528.75

# CUL row 91
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 7679.4
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 9224.25
$ws.Range("K91").Value = 4500
$ws.Range("L91").Value = 27672.75
$ws.Range("M91").Value = -3018
$ws.Range("N91").Value = -30636.75

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4266.242
$ws.Range("I107").Value = 368.13333
$ws.Range("J107").Value = 5412.745
$ws.Range("K107").Value = 1104.39999
$ws.Range("L107").Value = 16238.235
$ws.Range("M107").Value = 815.6000100000001
$ws.Range("N107").Value = -20078.235

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5360.0513
$ws.Range("I102").Value = 3472.0908
$ws.Range("J102").Value = 6101.75
$ws.Range("K102").Value = 3472.0908
$ws.Range("L102").Value = 6101.75
$ws.Range("M102").Value = -1850.0908
$ws.Range("N102").Value = -9345.75

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7640.4136
$ws.Range("I7").Value = 6609.5625
$ws.Range("J7").Value = 8909.154
$ws.Range("K7").Value = 6609.5625
$ws.Range("L7").Value = 8909.154
$ws.Range("M7").Value = -6497.5625
$ws.Range("N7").Value = -9133.154

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4491.2856
$ws.Range("I40").Value = 2839.4
$ws.Range("J40").Value = 8621
$ws.Range("K40").Value = 2839.4
$ws.Range("L40").Value = 8621
$ws.Range("M40").Value = -2703.4
$ws.Range("N40").Value = -8893

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6995.353
$ws.Range("I122").Value = 6196.6
$ws.Range("J122").Value = 7328.1665
$ws.Range("K122").Value = 18589.8
$ws.Range("L122").Value = 21984.4995
$ws.Range("M122").Value = -16139.8
$ws.Range("N122").Value = -26884.4995

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7640.4136
$ws.Range("I126").Value = 6609.5625
$ws.Range("J126").Value = 8909.154
$ws.Range("K126").Value = 19828.6875
$ws.Range("L126").Value = 26727.462
$ws.Range("M126").Value = -17358.6875
$ws.Range("N126").Value = -31667.462

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I136").Value = 35720548
$ws.Range("J136").Value = 7600
$ws.Range("K136").Value = 107161644
$ws.Range("L136").Value = 22800
$ws.Range("M136").Value = -107159094
$ws.Range("N136").Value = -27900

# LTW row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 59999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 59999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 59999
$ws.Range("N141").Value = -70359

# WVR row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 23942.143
$ws.Range("I26").Value = 19987
$ws.Range("J26").Value = 24601.334
$ws.Range("K26").Value = 19987
$ws.Range("L26").Value = 24601.334
$ws.Range("M26").Value = -19694
$ws.Range("N26").Value = -25187.334

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 755.4706
$ws.Range("I107").Value = 497.82608
$ws.Range("J107").Value = 1294.1818
$ws.Range("K107").Value = 1493.47824
$ws.Range("L107").Value = 3882.5454
$ws.Range("M107").Value = 426.5217600000001
$ws.Range("N107").Value = -7722.5454

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 11116408
$ws.Range("I113").Value = 20835178
$ws.Range("J113").Value = 9242.571
$ws.Range("K113").Value = 62505534
$ws.Range("L113").Value = 27727.713
$ws.Range("M113").Value = -62503364
$ws.Range("N113").Value = -32067.713

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19138.555
$ws.Range("I132").Value = 19350.8
$ws.Range("J132").Value = 18873.25
$ws.Range("K132").Value = 58052.39999999999
$ws.Range("L132").Value = 56619.75
$ws.Range("M132").Value = -55522.39999999999
$ws.Range("N132").Value = -61679.75
